$wb = $excel.ActiveWorkbook

# --- Sheet "main": bump Price (D2) to show two decimals ---
$wsMain = $wb.Worksheets.Item("main")
$wsMain.Range("D2").NumberFormat = "#,##0.00"

# --- Sheet "model": tidy up the percentage inputs to two decimals ---
$wsModel = $wb.Worksheets.Item("model")
$wsModel.Range("V17:V19").NumberFormat = "0.00%"
$wsModel.Range("V21").NumberFormat = "#,##0.00"

# --- Sheet "model": lower the ROIC assumption from 6% to 4% ---
$wsModel.Range("V17").Value = 0.04

# --- Restore selections / active cells to match where the edits were made ---
$wsMain.Range("D2").Select() | Out-Null
$wsModel.Activate() | Out-Null
$wsModel.Range("S21").Select() | Out-Null
